$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Starville Aingula" (row 5) was still carrying the sheet's default look;
# bring it in line with the other data rows (A2:D4) by copying their style.
$ws.Range("A5:D5").Style = $ws.Range("A2:D2").Style

# Add the new person reported below Starville Aingula.
$ws.Range("A6").Value = "Jake Opiyo"
$ws.Range("B6").Value = 29
$ws.Range("C6").Value = "Not Subscribed"
$ws.Range("D6").Value = "Employed"
